$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.018.23'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.43%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.263.34'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.59%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '252.81'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.33%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.631'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.07%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '71.15'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.48%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.677'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +16.62%  '
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.53'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.25%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0978'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.99%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '59.44'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.35%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.61'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.17%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.105'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.61%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.602.65'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.91%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.889'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.32%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.87'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.63%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.268.71'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.78%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '42.910.97'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.41%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0984'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.92%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.29'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.58%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.21'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.69%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.89'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.12'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.79%  '
$ws.Range("E25").Value = '  +1.36%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.78'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.08%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("E28").Value = '  -2.18%  '
$ws.Range("E29").Value = '  -1.23%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.15'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.80%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '167.96'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.22%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.25'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.20%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.30'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +14.35%  '
$ws.Range("E34").Value = '  +8.16%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0771'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.54%  '
$ws.Range("E36").Value = '  +1.66%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '28.92'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +11.30%  '
$ws.Range("E38").Value = '  +1.91%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.15'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.25%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0322'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.26%  '
$ws.Range("E41").Value = '  +3.35%  '
$ws.Range("E42").Value = '  +2.26%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.19'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.88%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '64.31'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.15%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.04'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.98%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.202'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.13%  '
$ws.Range("E47").Value = '  +1.90%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.103'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.98%  '
$ws.Range("E49").Value = '  -4.88%  '
$ws.Range("E50").Value = '  -0.03%  '

Write-Host "Updated cryptos list"